$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C12").Value = "CAPM_alpha_beta_test"
$ws.Range("A40").Value = "Bull/Bear beta"
$ws.Range("C40").Value = "Bull_Bear_beta_test"

$ws.Range("C12").Select()
